# Apply updated crypto price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to keep a literal text value even when the
    # string looks like a number (e.g. "561.35"), then drop the
    # temporary Text number-format so no stray style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '68.971.44'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '2.470.81'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  -0.02%  '
Set-TextValue $ws.Range("D5") '561.35'
$ws.Range("E5").Value = '  +0.97%  '
Set-TextValue $ws.Range("D6") '164.24'
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +2.42%  '
$ws.Range("D9").Value = '2.471.12'
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("E10").Value = '  +6.69%  '
$ws.Range("E11").Value = '  +0.98%  '
Set-TextValue $ws.Range("D12") '0.332'
$ws.Range("E12").Value = '  -0.77%  '
Set-TextValue $ws.Range("D13") '4.85'
$ws.Range("E13").Value = '  +2.36%  '
$ws.Range("D14").Value = '68.865.63'
$ws.Range("E14").Value = '  +0.94%  '
Set-TextValue $ws.Range("D15") '0.0000171'
$ws.Range("E15").Value = '  +2.91%  '
Set-TextValue $ws.Range("D16") '23.63'
$ws.Range("E16").Value = '  +2.58%  '
Set-TextValue $ws.Range("D17") '10.63'
$ws.Range("E17").Value = '  -1.56%  '
Set-TextValue $ws.Range("D18") '339.07'
$ws.Range("E18").Value = '  +0.25%  '
Set-TextValue $ws.Range("D19") '6.94'
$ws.Range("E19").Value = '  -1.91%  '
Set-TextValue $ws.Range("D20") '3.81'
$ws.Range("E20").Value = '  +2.22%  '
$ws.Range("E21").Value = '  +3.89%  '
$ws.Range("E22").Value = '  +0.05%  '
Set-TextValue $ws.Range("D23") '66.84'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("E24").Value = '  +1.69%  '
Set-TextValue $ws.Range("D25") '8.23'
$ws.Range("E25").Value = '  +2.81%  '
$ws.Range("D26").Value = '0.0₃0827'
$ws.Range("E26").Value = '  +1.16%  '
Set-TextValue $ws.Range("D27") '7.23'
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("E28").Value = '  +0.00%  '
Set-TextValue $ws.Range("D29") '430.16'
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("E30").Value = '  +0.60%  '
Set-TextValue $ws.Range("D31") '1.63'
$ws.Range("E31").Value = '  +0.44%  '
Set-TextValue $ws.Range("D32") '161.01'
$ws.Range("E32").Value = '  +2.46%  '
Set-TextValue $ws.Range("D33") '19.01'
$ws.Range("E33").Value = '  +0.12%  '
$ws.Range("E34").Value = '  +0.00%  '
Set-TextValue $ws.Range("D35") '0.107'
$ws.Range("E35").Value = '  -0.64%  '
Set-TextValue $ws.Range("D36") '17.90'
$ws.Range("E36").Value = '  +1.41%  '
Set-TextValue $ws.Range("D37") '4.43'
$ws.Range("E37").Value = '  +1.55%  '
$ws.Range("E38").Value = '  -0.64%  '
$ws.Range("E39").Value = '  +0.28%  '
Set-TextValue $ws.Range("D40") '1.07'
$ws.Range("E40").Value = '  +0.49%  '
Set-TextValue $ws.Range("D41") '2.06'
$ws.Range("E41").Value = '  +1.95%  '
Set-TextValue $ws.Range("D42") '3.38'
$ws.Range("E42").Value = '  +2.29%  '
Set-TextValue $ws.Range("D43") '131.10'
$ws.Range("E43").Value = '  -1.63%  '
$ws.Range("E44").Value = '  +1.44%  '
Set-TextValue $ws.Range("D45") '0.485'
$ws.Range("E45").Value = '  +2.20%  '
Set-TextValue $ws.Range("D46") '0.567'
$ws.Range("E46").Value = '  +1.70%  '
$ws.Range("E47").Value = '  +2.02%  '
$ws.Range("E48").Value = '  +0.35%  '
Set-TextValue $ws.Range("D49") '1.38'
$ws.Range("E50").Value = '  -3.50%  '
Set-TextValue $ws.Range("D51") '16.90'
$ws.Range("E51").Value = '  -1.45%  '
